$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recompute Price Impact (F), Incremental IL (G) and IL/Price Impact (H) using
# a forward-looking (next price tick) convention instead of the previous
# backward-looking one. The last simulated row (19) has no following tick,
# so its Price Impact / Incremental IL / ratio collapse to 0.

$ws.Range("F3").Value = -5.769230769230759
$ws.Range("G3").Value = -5.607440880921888
$ws.Range("H3").Value = 97.19564193597957

$ws.Range("F4").Value = -5.454545454545467
$ws.Range("G4").Value = -5.018312313404016
$ws.Range("H4").Value = 92.00239241240675

$ws.Range("F5").Value = -5.172413793103459
$ws.Range("G5").Value = -4.521339728624496
$ws.Range("H5").Value = 87.41256808674007

$ws.Range("F6").Value = -4.91803278688524
$ws.Range("G6").Value = -4.096946291451964
$ws.Range("H6").Value = 83.30457459285671

$ws.Range("F7").Value = -4.687499999999989
$ws.Range("G7").Value = -3.730659797764002
$ws.Range("H7").Value = 79.58740901896556

$ws.Range("F8").Value = -4.477611940298509
$ws.Range("G8").Value = -3.411562218035114
$ws.Range("H8").Value = 76.19155620278418

$ws.Range("F9").Value = -4.285714285714293
$ws.Range("G9").Value = -3.131269670282855
$ws.Range("H9").Value = 73.06295897326649

$ws.Range("F10").Value = -4.109589041095885
$ws.Range("G10").Value = -2.883243076158815
$ws.Range("H10").Value = 70.15891485319791

$ws.Range("F11").Value = -3.947368421052622
$ws.Range("G11").Value = -2.662311186069855
$ws.Range("H11").Value = 67.44521671376982

$ws.Range("F12").Value = -3.797468354430389
$ws.Range("G12").Value = -2.464333596800128
$ws.Range("H12").Value = 64.89411804906989

$ws.Range("F13").Value = -3.658536585365857
$ws.Range("G13").Value = -2.285958210903805
$ws.Range("H13").Value = 62.48285776470394

$ws.Range("F14").Value = -3.529411764705881
$ws.Range("G14").Value = -2.124443737774862
$ws.Range("H14").Value = 60.19257257028777

$ws.Range("F15").Value = -3.409090909090906
$ws.Range("G15").Value = -1.977527827510139
$ws.Range("H15").Value = 58.00748294029747

$ws.Range("F16").Value = -3.296703296703307
$ws.Range("G16").Value = -1.843327761696856
$ws.Range("H16").Value = 55.91427543813779

$ws.Range("F17").Value = -3.191489361702127
$ws.Range("G17").Value = -1.720264728470611
$ws.Range("H17").Value = 53.90162815874582

$ws.Range("F18").Value = -3.092783505154639
$ws.Range("G18").Value = -1.607005421057017
$ws.Range("H18").Value = 51.95984194751022

$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
